$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New results rows appended to the "fanduel entry history 20170120" sheet.
# Row 4-7: "NBA 50-50 5-Pack 2" style contest entries (50/50 Contest ($2 -
# Top 50% Win) (Main)), row 8: a separate beginners-only 50/50 contest.
# ---------------------------------------------------------------------------

# Row 4
$ws.Range("A4").Value = "S1147091687"
$ws.Range("B4").Value = "nba"
$ws.Range("C4").Value = 42755
$ws.Range("D4").Value = "50/50 Contest (`$2 - Top 50% Win) (Main)"
$ws.Range("E4").Value = "`$60k"
$ws.Range("F4").Value = 199.1
$ws.Range("H4").Value = 99
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = "Tournament"
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = "https://www.fanduel.com/entry/BMUMWPGJB"

# Row 5
$ws.Range("A5").Value = "S1147091218"
$ws.Range("B5").Value = "nba"
$ws.Range("C5").Value = 42755
$ws.Range("D5").Value = "50/50 Contest (`$2 - Top 50% Win) (Main)"
$ws.Range("E5").Value = "`$60k"
$ws.Range("F5").Value = 199.1
$ws.Range("H5").Value = 98
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = "Tournament"
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = "https://www.fanduel.com/entry/DAJIYWCAD"

# Row 6
$ws.Range("A6").Value = "S1147091098"
$ws.Range("B6").Value = "nba"
$ws.Range("C6").Value = 42755
$ws.Range("D6").Value = "50/50 Contest (`$2 - Top 50% Win) (Main)"
$ws.Range("E6").Value = "`$60k"
$ws.Range("F6").Value = 199.1
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = "Tournament"
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = "https://www.fanduel.com/entry/DJZPSNDRN"

# Row 7
$ws.Range("A7").Value = "S1147090959"
$ws.Range("B7").Value = "nba"
$ws.Range("C7").Value = 42755
$ws.Range("D7").Value = "50/50 Contest (`$2 - Top 50% Win) (Main)"
$ws.Range("E7").Value = "`$60k"
$ws.Range("F7").Value = 199.1
$ws.Range("H7").Value = 98
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = "Tournament"
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = "https://www.fanduel.com/entry/ADMWSKSTV"

# Row 8 is a brand-new row below the previous used range, so it has none of
# the existing column formatting (e.g. the date style on column C) yet --
# clone it from the row above before writing the value.
$null = $ws.Range("C7").Copy()
$null = $ws.Range("C8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A8").Value = "S1147089913"
$ws.Range("B8").Value = "nba"
$ws.Range("C8").Value = 42755
$ws.Range("D8").Value = "NBA 50/50 Contest (`$2_ Beginners Only) (Main)"
$ws.Range("E8").Value = "`$60k"
$ws.Range("F8").Value = 199.1
$ws.Range("H8").Value = 19
$ws.Range("I8").Value = 20
$ws.Range("J8").Value = "Tournament"
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = "https://www.fanduel.com/entry/ADMWOVNCL"

# Move the active selection to where the author last clicked.
$null = $ws.Range("L12").Select()
